$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.638665722778032
$ws.Range("K2").Value = 0.540250589217213
$ws.Range("L2").Value = 0.583644829216747
$ws.Range("N2").Value = 0.613804528949548

$ws.Range("B3").Value = 0.567431553913714
$ws.Range("F3").Value = 0.736611139000469
$ws.Range("K3").Value = 0.454763988800742
$ws.Range("L3").Value = 0.566244954822168
$ws.Range("N3").Value = 0.547449004875065

$ws.Range("B4").Value = 0.555344455602302
$ws.Range("K4").Value = 0.452271324998767
$ws.Range("L4").Value = 0.582472537227725
$ws.Range("N4").Value = 0.527045938028965

$ws.Range("B5").Value = 0.50130347410966
$ws.Range("K5").Value = 0.339691396434735
$ws.Range("L5").Value = 0.572223498223179
$ws.Range("N5").Value = 0.484225194294323

$ws.Range("B6").Value = 0.494042619036955
$ws.Range("K6").Value = 0.328236734725489
$ws.Range("L6").Value = 0.583489350682985
$ws.Range("N6").Value = 0.455937752846232

$ws.Range("B7").Value = 0.490655362928479
$ws.Range("K7").Value = 0.365701377881775
$ws.Range("L7").Value = 0.495615672214494
$ws.Range("N7").Value = 0.481729715314088

$ws.Range("B8").Value = 0.486367060267315
$ws.Range("J8").Value = 0.537044871026422
$ws.Range("K8").Value = 0.315893971140927
$ws.Range("L8").Value = 0.453289286617075
$ws.Range("N8").Value = 0.474756141069492

$ws.Range("B9").Value = 0.481241088566596
$ws.Range("K9").Value = 0.420810414510251
$ws.Range("L9").Value = 0.364456326664824
$ws.Range("N9").Value = 0.456964185261476

$ws.Range("B10").Value = 0.477954602429583
$ws.Range("K10").Value = 0.350511153488336
$ws.Range("L10").Value = 0.349016296615525
$ws.Range("N10").Value = 0.45006847466715

$ws.Range("B11").Value = 0.37503586986501
$ws.Range("K11").Value = 0.260747159004789
$ws.Range("L11").Value = 0.351533203159595
$ws.Range("N11").Value = 0.357192173089113
